# Weekly fruit/vegetable price update: add a new reporting date (45034)
# for "Betarraga" at "Terminal Hortofrutícola Agro Chillán", inserted as
# the two newest rows (Primera / Segunda quality), pushing the existing
# historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the first data row of this block
# (row 530), shifting the existing rows 530:601 down to 532:603.
$ws.Rows.Item(530).Insert()
$ws.Rows.Item(530).Insert()

# --- New row 530: Betarraga, Primera ---
$ws.Range("A530").Value = 7
$ws.Range("B530").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C530").Value = "Ñuble"
$ws.Range("D530").Value = 45034
$ws.Range("E530").Value = 16
$ws.Range("F530").Value = 100114014
$ws.Range("G530").Value = "Betarraga"
$ws.Range("H530").Value = "Sin especificar"
$ws.Range("I530").Value = "Primera"
$ws.Range("J530").Value = 300
$ws.Range("K530").Value = 1000
$ws.Range("L530").Value = 1000
$ws.Range("M530").Value = 1000
$ws.Range("N530").Value = "`$/paquete 5 unidades"
$ws.Range("O530").Value = "Provincia de Diguillín"
$ws.Range("P530").Value = 200
$ws.Range("Q530").Value = 5
$ws.Range("R530").Value = "Hortaliza"

# --- New row 531: Betarraga, Segunda ---
$ws.Range("A531").Value = 7
$ws.Range("B531").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C531").Value = "Ñuble"
$ws.Range("D531").Value = 45034
$ws.Range("E531").Value = 16
$ws.Range("F531").Value = 100114014
$ws.Range("G531").Value = "Betarraga"
$ws.Range("H531").Value = "Sin especificar"
$ws.Range("I531").Value = "Segunda"
$ws.Range("J531").Value = 300
$ws.Range("K531").Value = 800
$ws.Range("L531").Value = 800
$ws.Range("M531").Value = 800
$ws.Range("N531").Value = "`$/paquete 5 unidades"
$ws.Range("O531").Value = "Provincia de Diguillín"
$ws.Range("P531").Value = 160
$ws.Range("Q531").Value = 5
$ws.Range("R531").Value = "Hortaliza"
